# Generate Report for Archive
#
# This script:
#  1) Updates the localization status text "Ready for handoff" -> "In Translation"
#     on every sheet where it appears (Overview!E2:F2, zh-cn!C2, de-de!C2).
#  2) Shrinks the now-narrower "Status" columns to match the updated content
#     (Overview columns E & F, zh-cn column C, de-de column C).

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"
$newColWidth = 12.576851254417766

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = $newColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColWidth

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $newStatus
$wsZh.Columns.Item(3).ColumnWidth = $newColWidth

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $newStatus
$wsDe.Columns.Item(3).ColumnWidth = $newColWidth
